$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 315, shifting rows 315:346 down to 316:347
$ws.Rows.Item(315).Insert()

# Populate the new row 315 with the data from the diff
$ws.Range("A315").Value = 10
$ws.Range("B315").Value = "Vega Modelo de Temuco"
$ws.Range("C315").Value = "La Araucanía"
$ws.Range("D315").Value = 45166
$ws.Range("E315").Value = 9
$ws.Range("F315").Value = 100112013
$ws.Range("G315").Value = "Alcachofa"
$ws.Range("H315").Value = "Española"
$ws.Range("I315").Value = "Primera"
$ws.Range("J315").Value = 400
$ws.Range("K315").Value = 16000
$ws.Range("L315").Value = 16000
$ws.Range("M315").Value = 16000
$ws.Range("N315").Value = "$/caja 30 unidades"
$ws.Range("O315").Value = "Provincia de Limarí"
$ws.Range("P315").Value = 533
$ws.Range("Q315").Value = 30
$ws.Range("R315").Value = "Hortaliza"
